# "Testing - modelos con z-core"
#
# On the "full_signals - with decay" sheet (2nd sheet):
#   - the "alpha" column is removed entirely (old column H), shifting the
#     optimizer/metrics/tiempo columns one to the left
#   - the decay_steps formula for the existing model row now references the
#     max_epoch cell instead of a hard-coded literal
#   - two new model rows ("modelo 2 - norm z-core" and
#     "modelo 3 - norm z-core") are populated with their hyper-parameters /
#     metrics

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("full_signals - with decay")

# Drop the now-unused "alpha" column; everything to its right (optimizer,
# loss/nmse metrics, tiempo) shifts left by one column.
$ws.Columns("H:H").Delete()

# Existing "modelo 1" row: decay_steps now derives from max_epoch (E6)
# instead of the hard-coded 800.
$ws.Range("G6").Formula = "=80%*E6"

# New row 7: modelo 2 - norm z-core
$ws.Range("B7").Value = "modelo 2 - norm z-core"
$ws.Range("C7").Value = 30
$ws.Range("D7").Value = 0.0001
$ws.Range("E7").Value = 500
$ws.Range("F7").Value = 8
$ws.Range("G7").Formula = "=80%*E7"
$ws.Range("H7").Value = "Adam"
$ws.Range("J7").Value = 4.5789
$ws.Range("K7").Value = 0.256
$ws.Range("L7").Value = 4.1781
$ws.Range("M7").Value = 0.2122
$ws.Range("O7").Value = 120

# New row 8: modelo 3 - norm z-core
$ws.Range("B8").Value = "modelo 3 - norm z-core"
$ws.Range("C8").Value = 30
$ws.Range("D8").Value = 0.0001
$ws.Range("E8").Value = 410
$ws.Range("F8").Value = 8
$ws.Range("G8").Formula = "=80%*E8"
$ws.Range("H8").Value = "Adam"

$ws.Range("H8").Select() | Out-Null
